$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Overview sheet: row 3 (b.md) gets a fresh handoff "Ready for handoff"
# -----------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E3").Value = "Ready for handoff"
$ovw.Range("F3").Value = "Ready for handoff"
$ovw.Range("G3").Value = "2016-09-03 14:40:59"

# -----------------------------------------------------------------
# zh-cn sheet: row 3 (b.md) handoff refresh
# -----------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("H3").Value = "2016-09-03 14:40:55"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e9b89aa8d1275e6929011ecfdd89e77d53d869f6/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/237eedc497413a26d6a3388a4f2b01f7cc5dc936/e2e/b.md."
$zh.Columns.Item(16).ColumnWidth = 39.166666666666664

# -----------------------------------------------------------------
# de-de sheet: row 3 (b.md) handoff refresh
# -----------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("H3").Value = "2016-09-03 14:40:59"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e9b89aa8d1275e6929011ecfdd89e77d53d869f6/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/237eedc497413a26d6a3388a4f2b01f7cc5dc936/e2e/b.md."
$de.Columns.Item(16).ColumnWidth = 39.166666666666664
